$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfb1"
$ws.Range("C2").Value = "Itgb8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 88.72291666666666
$ws.Range("H2").Value = 266.16875
$ws.Range("I2").Value = 0.7675060578750151
$ws.Range("J2").Value = 0.7675060578750152
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1126243333333333
$ws.Range("N2").Value = 0.337873
$ws.Range("O2").Value = 0.01082936903163217
$ws.Range("P2").Value = 0.01082936903163217
$ws.Range("Q2").Value = 9.992359340972222
$ws.Range("R2").Value = 89.93123406874999
$ws.Range("S2").Value = 0.008311606334741775
$ws.Range("T2").Value = 0.008311606334741777

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfb1"
$ws.Range("C3").Value = "Itgb8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 88.72291666666666
$ws.Range("H3").Value = 266.16875
$ws.Range("I3").Value = 0.7675060578750151
$ws.Range("J3").Value = 0.7675060578750152
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.323421
$ws.Range("N3").Value = 9.970263
$ws.Range("O3").Value = 0.3195628457125252
$ws.Range("P3").Value = 0.3195628457125252
$ws.Range("Q3").Value = 294.8636044312499
$ws.Range("R3").Value = 2653.77243988125
$ws.Range("S3").Value = 0.2452664199561419
$ws.Range("T3").Value = 0.245266419956142

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tgfb1"
$ws.Range("C4").Value = "Itgb8"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 88.72291666666666
$ws.Range("H4").Value = 266.16875
$ws.Range("I4").Value = 0.7675060578750151
$ws.Range("J4").Value = 0.7675060578750152
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.963852666666667
$ws.Range("N4").Value = 20.891558
$ws.Range("O4").Value = 0.6696077852558425
$ws.Range("P4").Value = 0.6696077852558425
$ws.Range("Q4").Value = 617.8533198236111
$ws.Range("R4").Value = 5560.679878412499
$ws.Range("S4").Value = 0.5139280315841314
$ws.Range("T4").Value = 0.5139280315841315

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tgfb1"
$ws.Range("C5").Value = "Itgb8"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.91585
$ws.Range("H5").Value = 53.74755
$ws.Range("I5").Value = 0.1549827702197958
$ws.Range("J5").Value = 0.1549827702197958
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1126243333333333
$ws.Range("N5").Value = 0.337873
$ws.Range("O5").Value = 0.01082936903163217
$ws.Range("P5").Value = 0.01082936903163217
$ws.Range("Q5").Value = 2.01776066235
$ws.Range("R5").Value = 18.15984596115
$ws.Range("S5").Value = 0.001678365612254821
$ws.Range("T5").Value = 0.001678365612254821

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tgfb1"
$ws.Range("C6").Value = "Itgb8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 17.91585
$ws.Range("H6").Value = 53.74755
$ws.Range("I6").Value = 0.1549827702197958
$ws.Range("J6").Value = 0.1549827702197958
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.323421
$ws.Range("N6").Value = 9.970263
$ws.Range("O6").Value = 0.3195628457125252
$ws.Range("P6").Value = 0.3195628457125252
$ws.Range("Q6").Value = 59.54191212285
$ws.Range("R6").Value = 535.87720910565
$ws.Range("S6").Value = 0.04952673508784836
$ws.Range("T6").Value = 0.04952673508784837

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tgfb1"
$ws.Range("C7").Value = "Itgb8"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 17.91585
$ws.Range("H7").Value = 53.74755
$ws.Range("I7").Value = 0.1549827702197958
$ws.Range("J7").Value = 0.1549827702197958
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.963852666666667
$ws.Range("N7").Value = 20.891558
$ws.Range("O7").Value = 0.6696077852558425
$ws.Range("P7").Value = 0.6696077852558425
$ws.Range("Q7").Value = 124.7633397981
$ws.Range("R7").Value = 1122.8700581829
$ws.Range("S7").Value = 0.1037776695196926
$ws.Range("T7").Value = 0.1037776695196926

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tgfb1"
$ws.Range("C8").Value = "Itgb8"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.960212333333333
$ws.Range("H8").Value = 26.880637
$ws.Range("I8").Value = 0.077511171905189
$ws.Range("J8").Value = 0.07751117190518901
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1126243333333333
$ws.Range("N8").Value = 0.337873
$ws.Range("O8").Value = 0.01082936903163217
$ws.Range("P8").Value = 0.01082936903163217
$ws.Range("Q8").Value = 1.009137940566778
$ws.Range("R8").Value = 9.082241465101
$ws.Range("S8").Value = 0.000839397084635571
$ws.Range("T8").Value = 0.0008393970846355712

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tgfb1"
$ws.Range("C9").Value = "Itgb8"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.960212333333333
$ws.Range("H9").Value = 26.880637
$ws.Range("I9").Value = 0.077511171905189
$ws.Range("J9").Value = 0.07751117190518901
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.323421
$ws.Range("N9").Value = 9.970263
$ws.Range("O9").Value = 0.3195628457125252
$ws.Range("P9").Value = 0.3195628457125252
$ws.Range("Q9").Value = 29.778557833059
$ws.Range("R9").Value = 268.007020497531
$ws.Range("S9").Value = 0.02476969066853493
$ws.Range("T9").Value = 0.02476969066853494

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tgfb1"
$ws.Range("C10").Value = "Itgb8"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.960212333333333
$ws.Range("H10").Value = 26.880637
$ws.Range("I10").Value = 0.077511171905189
$ws.Range("J10").Value = 0.07751117190518901
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.963852666666667
$ws.Range("N10").Value = 20.891558
$ws.Range("O10").Value = 0.6696077852558425
$ws.Range("P10").Value = 0.6696077852558425
$ws.Range("Q10").Value = 62.39759855138288
$ws.Range("R10").Value = 561.578386962446
$ws.Range("S10").Value = 0.05190208415201849
$ws.Range("T10").Value = 0.0519020841520185
